# Refresh the daily "cryptos" price/volume snapshot.
# Prices in column D are stored as text (they use '.' as both thousands
# and decimal separators), so numeric-looking values are written with a
# leading apostrophe to force Excel to keep them as literal text instead
# of silently reinterpreting/reformatting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'" + '38.862.47'
$ws.Range('E2').Value = '  +2.84%  '

# Row 3
$ws.Range('D3').Value = "'" + '2.091.87'
$ws.Range('E3').Value = '  +2.22%  '

# Row 4
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('D5').Value = "'" + '228.70'
$ws.Range('E5').Value = '  +0.38%  '

# Row 6
$ws.Range('D6').Value = "'" + '0.614'
$ws.Range('E6').Value = '  +0.56%  '

# Row 7
$ws.Range('D7').Value = "'" + '60.32'
$ws.Range('E7').Value = '  +0.81%  '

# Row 8
$ws.Range('E8').Value = '  -0.04%  '

# Row 9
$ws.Range('E9').Value = '  +1.87%  '

# Row 10
$ws.Range('D10').Value = "'" + '0.0841'
$ws.Range('E10').Value = '  +0.61%  '

# Row 11
$ws.Range('D11').Value = "'" + '0.104'
$ws.Range('E11').Value = '  -0.13%  '

# Row 12
$ws.Range('D12').Value = "'" + '2.401.27'
$ws.Range('E12').Value = '  +2.11%  '

# Row 13
$ws.Range('D13').Value = "'" + '15.00'
$ws.Range('E13').Value = '  +4.16%  '

# Row 14
$ws.Range('D14').Value = "'" + '21.97'
$ws.Range('E14').Value = '  +2.42%  '

# Row 15
$ws.Range('E15').Value = '  +4.23%  '

# Row 16
$ws.Range('E16').Value = '  -0.52%  '

# Row 17
$ws.Range('D17').Value = "'" + '2.089.67'
$ws.Range('E17').Value = '  +1.12%  '

# Row 18
$ws.Range('D18').Value = "'" + '38.776.39'
$ws.Range('E18').Value = '  +2.58%  '

# Row 19
$ws.Range('D19').Value = "'" + '71.59'
$ws.Range('E19').Value = '  +3.02%  '

# Row 20
$ws.Range('E20').Value = '  +2.26%  '

# Row 21
$ws.Range('D21').Value = '0.0₃0839'
$ws.Range('E21').Value = '  +1.03%  '

# Row 22
$ws.Range('D22').Value = "'" + '227.26'
$ws.Range('E22').Value = '  +2.19%  '

# Row 23
$ws.Range('E23').Value = '  -0.42%  '

# Row 24
$ws.Range('E24').Value = '  -0.21%  '

# Row 25
$ws.Range('E25').Value = '  +2.49%  '

# Row 26
$ws.Range('D26').Value = "'" + '171.07'
$ws.Range('E26').Value = '  +1.41%  '

# Row 27
$ws.Range('D27').Value = "'" + '9.54'
$ws.Range('E27').Value = '  +2.30%  '

# Row 28
$ws.Range('E28').Value = '  +10.39%  '

# Row 29
$ws.Range('D29').Value = "'" + '1.48'
$ws.Range('E29').Value = '  +14.20%  '

# Row 30
$ws.Range('E30').Value = '  +2.05%  '

# Row 31
$ws.Range('D31').Value = "'" + '0.121'
$ws.Range('E31').Value = '  +1.11%  '

# Row 32
$ws.Range('E32').Value = '  +5.37%  '

# Row 33
$ws.Range('E33').Value = '  +2.71%  '

# Row 34
$ws.Range('E34').Value = '  +3.87%  '

# Row 35
$ws.Range('E35').Value = '  +1.48%  '

# Row 36
$ws.Range('D36').Value = "'" + '6.49'
$ws.Range('E36').Value = '  +0.10%  '

# Row 37
$ws.Range('E37').Value = '  +1.12%  '

# Row 38
$ws.Range('D38').Value = "'" + '3.59'
$ws.Range('E38').Value = '  +2.91%  '

# Row 39
$ws.Range('E39').Value = '  -0.14%  '

# Row 40
$ws.Range('D40').Value = "'" + '18.14'
$ws.Range('E40').Value = '  -1.05%  '

# Row 41
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = "'" + '1.543.51'
$ws.Range('E41').Value = '  +1.04%  '

# Row 42
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = "'" + '0.0226'
$ws.Range('E42').Value = '  +4.59%  '

# Row 43
$ws.Range('D43').Value = "'" + '100.90'
$ws.Range('E43').Value = '  +3.25%  '

# Row 44
$ws.Range('E44').Value = '  -0.84%  '

# Row 45
$ws.Range('D45').Value = "'" + '0.0921'
$ws.Range('E45').Value = '  +3.52%  '

# Row 46
$ws.Range('D46').Value = "'" + '7.69'
$ws.Range('E46').Value = '  +8.19%  '

# Row 47
$ws.Range('D47').Value = "'" + '1.13'
$ws.Range('E47').Value = '  +1.76%  '

# Row 48
$ws.Range('D48').Value = "'" + '4.10'
$ws.Range('E48').Value = '  -1.78%  '

# Row 49
$ws.Range('E49').Value = '  +2.87%  '

# Row 50
$ws.Range('D50').Value = "'" + '2.97'
$ws.Range('E50').Value = '  +0.77%  '

# Row 51
$ws.Range('D51').Value = "'" + '2.288.64'
$ws.Range('E51').Value = '  +2.15%  '

